$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must remain text,
# matching the original inline-string (non-numeric) cell format.
$textForceCells = @("D5", "D8", "D10", "D11", "D14", "D15", "D17", "D18", "D20", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D39", "D40", "D42", "D44", "D45", "D46", "D48", "D49", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.938.91"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.555.82"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "206.90"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "21.98"
$ws.Range("E8").Value = "  +2.20%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "0.0588"
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").Value = "0.0859"
$ws.Range("D12").Value = "1.777.63"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "1.555.88"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "3.74"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").Value = "0.519"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "26.935.48"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "61.82"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "218.54"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "0.0₃0696"
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("D20").Value = "7.31"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "4.06"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").Value = "9.20"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").Value = "154.07"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").Value = "6.61"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").Value = "14.98"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "0.104"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "0.0470"
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("D31").Value = "1.10"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").Value = "1.438.91"
$ws.Range("E33").Value = "  +5.12%  "
$ws.Range("D34").Value = "3.07"
$ws.Range("E34").Value = "  +4.32%  "
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  +3.46%  "
$ws.Range("D36").Value = "0.979"
$ws.Range("E36").Value = "  +2.35%  "
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "0.522"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Value = "0.813"
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "5.68"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("D44").Value = "0.986"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "64.37"
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("D46").Value = "1.74"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("D47").Value = "1.691.07"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").Value = "87.32"
$ws.Range("E48").Value = "  +2.70%  "
$ws.Range("D49").Value = "0.0520"
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("D50").Value = "0.0₇0998"
$ws.Range("E50").Value = "  +3.13%  "
$ws.Range("D51").Value = "0.0961"
$ws.Range("E51").Value = "  +1.32%  "

# Restore default (un-styled) appearance now that the text is committed,
# so no stray number-format style lingers on these cells.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
